$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that look like plain numbers need to be forced to text storage
# (NumberFormat "@") so Excel keeps them as shared-string literals instead
# of silently converting to numeric cells; Style is then reset back to
# "Normal" so no stray style index is left on the cell.
$numericLooking = @("B2", "B3", "B4", "C4", "D2")
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("B2").Value = "0.17"
$ws.Range("B3").Value = "-0.01"
$ws.Range("B4").Value = "-0.09"
$ws.Range("C2").Value = "44.29***"
$ws.Range("C3").Value = "2.21***"
$ws.Range("C4").Value = "0.98"
$ws.Range("D2").Value = "-0.89"
$ws.Range("D3").Value = "0.46***"
$ws.Range("D4").Value = "0.82*"

foreach ($addr in $numericLooking) {
    $ws.Range($addr).Style = "Normal"
}
